$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(126).Insert()

$ws.Range("A126").Value = 5
$ws.Range("B126").Value = "Macroferia Regional de Talca"
$ws.Range("C126").Value = "Maule"
$ws.Range("D126").Value = 44767
$ws.Range("E126").Value = 7
$ws.Range("F126").Value = 100112017
$ws.Range("G126").Value = "Apio"
$ws.Range("H126").Value = "Americana (o)"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = 12000
$ws.Range("N126").Value = "`$/docena de matas"
$ws.Range("O126").Value = "Provincia del Elquí"
$ws.Range("P126").Value = 2000
$ws.Range("Q126").Value = 6
$ws.Range("R126").Value = "Hortaliza"
